$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 7.714752999999999
$ws.Cells.Item(2, 8).Value = 23.144259
$ws.Cells.Item(2, 9).Value = 0.01072102934448851
$ws.Cells.Item(2, 10).Value = 0.01072102934448851
$ws.Cells.Item(2, 13).Value = 1.819857
$ws.Cells.Item(2, 14).Value = 5.459571
$ws.Cells.Item(2, 15).Value = 0.01485317462584607
$ws.Cells.Item(2, 16).Value = 0.01485317462584607
$ws.Cells.Item(2, 17).Value = 14.039747250321
$ws.Cells.Item(2, 18).Value = 126.357725252889
$ws.Cells.Item(2, 19).Value = 0.0001592413210225079
$ws.Cells.Item(2, 20).Value = 0.0001592413210225079

# Row 3
$ws.Cells.Item(3, 7).Value = 7.714752999999999
$ws.Cells.Item(3, 8).Value = 23.144259
$ws.Cells.Item(3, 9).Value = 0.01072102934448851
$ws.Cells.Item(3, 10).Value = 0.01072102934448851
$ws.Cells.Item(3, 15).Value = 0.726618572334523
$ws.Cells.Item(3, 16).Value = 0.7266185723345231
$ws.Cells.Item(3, 17).Value = 686.8256355926799
$ws.Cells.Item(3, 18).Value = 6181.43072033412
$ws.Cells.Item(3, 19).Value = 0.00779009903624877
$ws.Cells.Item(3, 20).Value = 0.00779009903624877

# Row 4
$ws.Cells.Item(4, 7).Value = 7.714752999999999
$ws.Cells.Item(4, 8).Value = 23.144259
$ws.Cells.Item(4, 9).Value = 0.01072102934448851
$ws.Cells.Item(4, 10).Value = 0.01072102934448851
$ws.Cells.Item(4, 13).Value = 31.52924033333333
$ws.Cells.Item(4, 14).Value = 94.58772099999999
$ws.Cells.Item(4, 15).Value = 0.257333028084772
$ws.Cells.Item(4, 16).Value = 0.257333028084772
$ws.Cells.Item(4, 17).Value = 243.2403014493043
$ws.Cells.Item(4, 18).Value = 2189.162713043739
$ws.Cells.Item(4, 19).Value = 0.002758874945402928
$ws.Cells.Item(4, 20).Value = 0.002758874945402927

# Row 5
$ws.Cells.Item(5, 7).Value = 7.714752999999999
$ws.Cells.Item(5, 8).Value = 23.144259
$ws.Cells.Item(5, 9).Value = 0.01072102934448851
$ws.Cells.Item(5, 10).Value = 0.01072102934448851
$ws.Cells.Item(5, 11).Value = 1.0
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.1464426666666667
$ws.Cells.Item(5, 14).Value = 0.439328
$ws.Cells.Item(5, 15).Value = 0.001195224954858853
$ws.Cells.Item(5, 16).Value = 0.001195224954858853
$ws.Cells.Item(5, 17).Value = 1.129769001994666
$ws.Cells.Item(5, 18).Value = 10.167921017952
$ws.Cells.Item(5, 19).Value = 0.00001281404181430671
$ws.Cells.Item(5, 20).Value = 0.00001281404181430671

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9638361163234916
$ws.Cells.Item(6, 10).Value = 0.9638361163234914
$ws.Cells.Item(6, 13).Value = 1.819857
$ws.Cells.Item(6, 14).Value = 5.459571
$ws.Cells.Item(6, 15).Value = 0.01485317462584607
$ws.Cells.Item(6, 16).Value = 0.01485317462584607
$ws.Cells.Item(6, 17).Value = 1262.193678340166
$ws.Cells.Item(6, 18).Value = 11359.74310506149
$ws.Cells.Item(6, 19).Value = 0.0143160261464501
$ws.Cells.Item(6, 20).Value = 0.0143160261464501

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9638361163234916
$ws.Cells.Item(7, 10).Value = 0.9638361163234914
$ws.Cells.Item(7, 15).Value = 0.726618572334523
$ws.Cells.Item(7, 16).Value = 0.7266185723345231
$ws.Cells.Item(7, 19).Value = 0.7003412228074267
$ws.Cells.Item(7, 20).Value = 0.7003412228074266

# Row 8
$ws.Cells.Item(8, 9).Value = 0.9638361163234916
$ws.Cells.Item(8, 10).Value = 0.9638361163234914
$ws.Cells.Item(8, 13).Value = 31.52924033333333
$ws.Cells.Item(8, 14).Value = 94.58772099999999
$ws.Cells.Item(8, 15).Value = 0.257333028084772
$ws.Cells.Item(8, 16).Value = 0.257333028084772
$ws.Cells.Item(8, 17).Value = 21867.65654202562
$ws.Cells.Item(8, 18).Value = 196808.9088782305
$ws.Cells.Item(8, 19).Value = 0.2480268663909907
$ws.Cells.Item(8, 20).Value = 0.2480268663909906

# Row 9
$ws.Cells.Item(9, 9).Value = 0.9638361163234916
$ws.Cells.Item(9, 10).Value = 0.9638361163234914
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.1464426666666667
$ws.Cells.Item(9, 14).Value = 0.439328
$ws.Cells.Item(9, 15).Value = 0.001195224954858853
$ws.Cells.Item(9, 16).Value = 0.001195224954858853
$ws.Cells.Item(9, 17).Value = 101.5678748967324
$ws.Cells.Item(9, 18).Value = 914.1108740705919
$ws.Cells.Item(9, 19).Value = 0.001152000978624077
$ws.Cells.Item(9, 20).Value = 0.001152000978624077

# Row 10
$ws.Cells.Item(10, 7).Value = 15.35226833333333
$ws.Cells.Item(10, 8).Value = 46.056805
$ws.Cells.Item(10, 9).Value = 0.02133472313450974
$ws.Cells.Item(10, 10).Value = 0.02133472313450974
$ws.Cells.Item(10, 13).Value = 1.819857
$ws.Cells.Item(10, 14).Value = 5.459571
$ws.Cells.Item(10, 15).Value = 0.01485317462584607
$ws.Cells.Item(10, 16).Value = 0.01485317462584607
$ws.Cells.Item(10, 17).Value = 27.938932992295
$ws.Cells.Item(10, 18).Value = 251.450396930655
$ws.Cells.Item(10, 19).Value = 0.0003168883683109512
$ws.Cells.Item(10, 20).Value = 0.0003168883683109512

# Row 11
$ws.Cells.Item(11, 7).Value = 15.35226833333333
$ws.Cells.Item(11, 8).Value = 46.056805
$ws.Cells.Item(11, 9).Value = 0.02133472313450974
$ws.Cells.Item(11, 10).Value = 0.02133472313450974
$ws.Cells.Item(11, 15).Value = 0.726618572334523
$ws.Cells.Item(11, 16).Value = 0.7266185723345231
$ws.Cells.Item(11, 17).Value = 1366.774990181933
$ws.Cells.Item(11, 18).Value = 12300.9749116374
$ws.Cells.Item(11, 19).Value = 0.01550220606514979
$ws.Cells.Item(11, 20).Value = 0.01550220606514979

# Row 12
$ws.Cells.Item(12, 7).Value = 15.35226833333333
$ws.Cells.Item(12, 8).Value = 46.056805
$ws.Cells.Item(12, 9).Value = 0.02133472313450974
$ws.Cells.Item(12, 10).Value = 0.02133472313450974
$ws.Cells.Item(12, 13).Value = 31.52924033333333
$ws.Cells.Item(12, 14).Value = 94.58772099999999
$ws.Cells.Item(12, 15).Value = 0.257333028084772
$ws.Cells.Item(12, 16).Value = 0.257333028084772
$ws.Cells.Item(12, 17).Value = 484.0453579434894
$ws.Cells.Item(12, 18).Value = 4356.408221491405
$ws.Cells.Item(12, 19).Value = 0.005490128907553632
$ws.Cells.Item(12, 20).Value = 0.005490128907553631

# Row 13
$ws.Cells.Item(13, 7).Value = 15.35226833333333
$ws.Cells.Item(13, 8).Value = 46.056805
$ws.Cells.Item(13, 9).Value = 0.02133472313450974
$ws.Cells.Item(13, 10).Value = 0.02133472313450974
$ws.Cells.Item(13, 11).Value = 1.0
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.1464426666666667
$ws.Cells.Item(13, 14).Value = 0.439328
$ws.Cells.Item(13, 15).Value = 0.001195224954858853
$ws.Cells.Item(13, 16).Value = 0.001195224954858853
$ws.Cells.Item(13, 17).Value = 2.248227114115556
$ws.Cells.Item(13, 18).Value = 20.23404402704
$ws.Cells.Item(13, 19).Value = 0.00002549979349537053
$ws.Cells.Item(13, 20).Value = 0.00002549979349537053

# Row 14
$ws.Cells.Item(14, 7).Value = 2.956173
$ws.Cells.Item(14, 8).Value = 8.868519000000001
$ws.Cells.Item(14, 9).Value = 0.004108131197510101
$ws.Cells.Item(14, 10).Value = 0.0041081311975101
$ws.Cells.Item(14, 13).Value = 1.819857
$ws.Cells.Item(14, 14).Value = 5.459571
$ws.Cells.Item(14, 15).Value = 0.01485317462584607
$ws.Cells.Item(14, 16).Value = 0.01485317462584607
$ws.Cells.Item(14, 17).Value = 5.379812127261
$ws.Cells.Item(14, 18).Value = 48.41830914534901
$ws.Cells.Item(14, 19).Value = 0.00006101879006250365
$ws.Cells.Item(14, 20).Value = 0.00006101879006250364

# Row 15
$ws.Cells.Item(15, 7).Value = 2.956173
$ws.Cells.Item(15, 8).Value = 8.868519000000001
$ws.Cells.Item(15, 9).Value = 0.004108131197510101
$ws.Cells.Item(15, 10).Value = 0.0041081311975101
$ws.Cells.Item(15, 15).Value = 0.726618572334523
$ws.Cells.Item(15, 16).Value = 0.7266185723345231
$ws.Cells.Item(15, 17).Value = 263.18086912788
$ws.Cells.Item(15, 18).Value = 2368.62782215092
$ws.Cells.Item(15, 19).Value = 0.002985044425697704
$ws.Cells.Item(15, 20).Value = 0.002985044425697704

# Row 16
$ws.Cells.Item(16, 7).Value = 2.956173
$ws.Cells.Item(16, 8).Value = 8.868519000000001
$ws.Cells.Item(16, 9).Value = 0.004108131197510101
$ws.Cells.Item(16, 10).Value = 0.0041081311975101
$ws.Cells.Item(16, 13).Value = 31.52924033333333
$ws.Cells.Item(16, 14).Value = 94.58772099999999
$ws.Cells.Item(16, 15).Value = 0.257333028084772
$ws.Cells.Item(16, 16).Value = 0.257333028084772
$ws.Cells.Item(16, 17).Value = 93.205888983911
$ws.Cells.Item(16, 18).Value = 838.853000855199
$ws.Cells.Item(16, 19).Value = 0.001057157840824795
$ws.Cells.Item(16, 20).Value = 0.001057157840824795

# Row 17
$ws.Cells.Item(17, 7).Value = 2.956173
$ws.Cells.Item(17, 8).Value = 8.868519000000001
$ws.Cells.Item(17, 9).Value = 0.004108131197510101
$ws.Cells.Item(17, 10).Value = 0.0041081311975101
$ws.Cells.Item(17, 11).Value = 1.0
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.1464426666666667
$ws.Cells.Item(17, 14).Value = 0.439328
$ws.Cells.Item(17, 15).Value = 0.001195224954858853
$ws.Cells.Item(17, 16).Value = 0.001195224954858853
$ws.Cells.Item(17, 17).Value = 0.432909857248
$ws.Cells.Item(17, 18).Value = 3.896188715232
$ws.Cells.Item(17, 19).Value = 0.000004910140925098255
$ws.Cells.Item(17, 20).Value = 0.000004910140925098254
